# Swedish_data.xlsx — "Add files via upload / current versions of
# structural tables, to be used during non-presentation".
#
# The "comms" column (M) values get moved into the "comms_internal" column
# (S) wherever "comms" held a value, clearing the source cell. The
# "comms_internal" column also gets its direct formatting normalised to the
# one style already used by its neighbours. Finally, the saved selection
# moves from M19 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Normalize the format of the whole "comms_internal" column
#        (S2:S131) onto the single direct style already used by the rest of
#        the data rows (copy it from a neighbouring cell, e.g. N2, so the
#        same style entry is reused instead of minting new ones).
$ws.Range("N2").Copy()
$ws.Range("S2:S131").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Move each non-empty "comms" (column M) value over to
#        "comms_internal" (column S), clearing the original cell.
$rowsWithComms = @(
    2,3,4,6,12,13,14,15,20,22,23,25,30,33,34,38,43,44,46,47,49,51,52,57,59,
    62,65,66,69,77,80,81,83,84,88,90,91,93,95,97,98,99,100,103,104,105,107,
    112,113,114,115,116,117,118,119,120,121,123,124,125,126,128,129,130,131
)

foreach ($r in $rowsWithComms) {
    $mCell = $ws.Cells.Item($r, 13)   # column M = comms
    $sCell = $ws.Cells.Item($r, 19)   # column S = comms_internal
    $sCell.Value2 = $mCell.Value2
    $mCell.Value2 = $null
}

# --- 3. Update the saved view state: select A2 (was M19).
$ws.Range("A2").Select()
